$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rename business* labels to shorter names ---
$ws.Range("A1").Value = "email"
$ws.Range("C1").Value = "website"
$ws.Range("D1").Value = "phoneNumber"
$ws.Range("E1").Value = "address"
$ws.Range("F1").Value = "socialMedia"

# --- Restyle the F1 (socialMedia) header font (Courier New, ~9.8pt, #C77DBB) ---
$ws.Range("F1").Font.Name = "Courier New"
$ws.Range("F1").Font.Size = 9.8
$ws.Range("F1").Font.Color = 12287431

# --- Remove the pre-existing hyperlink on A2 before re-entering row data ---
$ws.Range("A2").Hyperlinks.Delete()

# --- Row 2: contact entry (now with a social handle + website hyperlink) ---
$ws.Range("A2").Value = "samplhghe@sa.com"
$ws.Range("B2").Value = "ryry"
$ws.Range("C2").Value = "www.das.com"
$ws.Range("D2").Value = 99293992
$ws.Range("E2").Value = 400101
$ws.Range("F2").Value = "-"

# --- Row 3: new contact entry ---
$ws.Range("A3").Value = "aqwe@dsa.com"
$ws.Range("B3").Value = "dsdf"
$ws.Range("C3").Value = "www.sds.com"
$ws.Range("D3").Value = 99293992
$ws.Range("E3").Value = 400101
$ws.Range("F3").Value = "sample"

# --- Row 4: new sample entry (re-uses the original A2 hyperlink + display) ---
# The hyperlink carries a stale "wfsd2@gmail.com" display/tooltip from the
# link it was copied from, but the cell text itself is "sample2" - so add
# the hyperlink (which seeds that stale display) before overwriting the
# visible cell value.
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:wfsd2@gmail.com", [Type]::Missing, [Type]::Missing, "wfsd2@gmail.com")
$ws.Range("A4").Value = "sample2"
$ws.Range("B4").Value = "sample"
$ws.Range("C4").Value = "sample"
$ws.Range("D4").Value = 99293992
$ws.Range("E4").Value = 400101
$ws.Range("F4").Value = "sample"

# --- Row 5: new sample entry ---
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:wfsd2@gmail.com", [Type]::Missing, [Type]::Missing, "wfsd2@gmail.com")
$ws.Range("A5").Value = "sample3"
$ws.Range("B5").Value = "sample"
$ws.Range("C5").Value = "sample"
$ws.Range("D5").Value = 99293992
$ws.Range("E5").Value = 400101
$ws.Range("F5").Value = "sample"

# --- Row 6: new sample entry ---
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:wfsd2@gmail.com", [Type]::Missing, [Type]::Missing, "wfsd2@gmail.com")
$ws.Range("A6").Value = "sample4"
$ws.Range("B6").Value = "sample"
$ws.Range("C6").Value = "sample"
$ws.Range("D6").Value = 99293992
$ws.Range("E6").Value = 400101
$ws.Range("F6").Value = "sample"

# --- Remaining hyperlinks (website + email columns, no stale display text) ---
$ws.Hyperlinks.Add($ws.Range("C2"), "www.das.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "www.sds.com")

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:samplhghe@sa.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:aqwe@dsa.com")

# Re-apply the shared "Hyperlink" cell style everywhere Hyperlinks.Add just
# minted its own duplicate style entry (matches the pre-existing look used
# across column A).
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("A4").Style = "Hyperlink"
$ws.Range("A5").Style = "Hyperlink"
$ws.Range("A6").Style = "Hyperlink"
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("C3").Style = "Hyperlink"

# --- Selection moves to C2 ---
$ws.Range("C2").Select()
